$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 253, shifting existing rows 253-274 down to 254-275.
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row 253 with the new data record.
$ws.Cells.Item(253, 1).Value = 5
$ws.Cells.Item(253, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(253, 3).Value = "Maule"
$ws.Cells.Item(253, 4).Value = 44783
$ws.Cells.Item(253, 5).Value = 7
$ws.Cells.Item(253, 6).Value = "Fruta"
$ws.Cells.Item(253, 7).Value = 100108
$ws.Cells.Item(253, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(253, 9).Value = 100108005
$ws.Cells.Item(253, 10).Value = "Piña"
$ws.Cells.Item(253, 11).Value = "Caramelo"
$ws.Cells.Item(253, 12).Value = "Tercera"
$ws.Cells.Item(253, 13).Value = 250
$ws.Cells.Item(253, 14).Value = 19000
$ws.Cells.Item(253, 15).Value = 19000
$ws.Cells.Item(253, 16).Value = 19000
$ws.Cells.Item(253, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(253, 18).Value = "Ecuador"
$ws.Cells.Item(253, 19).Value = 1188
$ws.Cells.Item(253, 20).Value = 16

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(253, 4).NumberFormat = $ws.Cells.Item(254, 4).NumberFormat
